# Insert a new data row for "Poroto verde" (Macroferia Regional de Talca)
# above the current row 127. This shifts all subsequent rows (old rows
# 127-217) down by one (to 128-218), which is exactly the pattern shown
# in the diff, and extends the sheet's used range to A1:R218.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 127 (and everything below it) down by one row.
$ws.Rows("127:127").Insert()

# Populate the newly inserted row 127 with the new record.
$ws.Range("A127").Value = 5
$ws.Range("B127").Value = "Macroferia Regional de Talca"
$ws.Range("C127").Value = "Maule"
$ws.Range("D127").Value = 44957
$ws.Range("E127").Value = 7
$ws.Range("F127").Value = 100112031
$ws.Range("G127").Value = "Poroto verde"
$ws.Range("H127").Value = "Sin especificar"
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 150
$ws.Range("K127").Value = 25000
$ws.Range("L127").Value = 25000
$ws.Range("M127").Value = 25000
$ws.Range("N127").Value = "`$/saco 25 kilos"
$ws.Range("O127").Value = "Región del Maule"
$ws.Range("P127").Value = 1000
$ws.Range("Q127").Value = 25
$ws.Range("R127").Value = "Hortaliza"
